# Generate Report for Handoff
#
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# - Refreshed timestamps:
#     Overview!G2 (Latest HO Xliff Generate Date)      2016-09-01 01:03:11 -> 2016-09-01 01:03:57
#     zh-cn!H2    (Latest Handoff Datetime)             2016-09-01 01:03:04 -> 2016-09-01 01:03:53
# - Columns that held the long status text are narrower now that the text is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"

# --- Status text updates ---
$wsOverview.Range("E2").Value = $readyForHandoff
$wsOverview.Range("F2").Value = $readyForHandoff
$wsZhCn.Range("C2").Value     = $readyForHandoff
$wsDeDe.Range("C2").Value     = $readyForHandoff

# --- Timestamp updates (stored as text, not real dates) ---
$wsOverview.Range("G2").Value = "2016-09-01 01:03:57"
$wsZhCn.Range("H2").Value     = "2016-09-01 01:03:53"

# --- Column width refresh, now that the status text is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3333333333333
